$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.420.08"
$ws.Range("E2").Value = "  +0.64%  "
$ws.Range("D3").Value = "2.267.22"
$ws.Range("E3").Value = "  -0.50%  "
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").Value = "'120.70"
$ws.Range("E5").Value = "  +7.36%  "
$ws.Range("D6").Value = "'268.85"
$ws.Range("E6").Value = "  +1.25%  "
$ws.Range("D7").Value = "'0.641"
$ws.Range("E7").Value = "  +3.07%  "
$ws.Range("D8").Value = "'1.01"
$ws.Range("E8").Value = "  +0.19%  "
$ws.Range("D9").Value = "'0.620"
$ws.Range("E9").Value = "  +2.04%  "
$ws.Range("D10").Value = "'47.54"
$ws.Range("E10").Value = "  -0.29%  "
$ws.Range("D11").Value = "'0.0942"
$ws.Range("E11").Value = "  +1.26%  "
$ws.Range("D12").Value = "'9.39"
$ws.Range("E12").Value = "  +5.63%  "
$ws.Range("E13").Value = "  -1.80%  "
$ws.Range("D14").Value = "'15.80"
$ws.Range("E14").Value = "  +1.80%  "
$ws.Range("D15").Value = "'0.912"
$ws.Range("E15").Value = "  +6.36%  "
$ws.Range("D16").Value = "2.610.99"
$ws.Range("E16").Value = "  -0.43%  "
$ws.Range("D17").Value = "2.265.06"
$ws.Range("E17").Value = "  -0.56%  "
$ws.Range("D18").Value = "43.593.54"
$ws.Range("E18").Value = "  +0.78%  "
$ws.Range("D19").Value = "'0.0000110"
$ws.Range("E19").Value = "  +1.41%  "
$ws.Range("D20").Value = "'6.92"
$ws.Range("E20").Value = "  +1.93%  "
$ws.Range("D21").Value = "'72.70"
$ws.Range("E21").Value = "  +1.70%  "
$ws.Range("E22").Value = "  -5.44%  "
$ws.Range("D23").Value = "'234.90"
$ws.Range("E23").Value = "  +1.24%  "
$ws.Range("D24").Value = "'2.95"
$ws.Range("E24").Value = "  +3.14%  "
$ws.Range("D25").Value = "'9.62"
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("D26").Value = "'12.22"
$ws.Range("E26").Value = "  +7.93%  "
$ws.Range("E27").Value = "  +0.93%  "
$ws.Range("D28").Value = "'42.19"
$ws.Range("E28").Value = "  +4.36%  "
$ws.Range("D29").Value = "'3.35"
$ws.Range("E29").Value = "  +0.15%  "
$ws.Range("E30").Value = "  +0.30%  "
$ws.Range("D31").Value = "'175.06"
$ws.Range("E31").Value = "  +1.65%  "
$ws.Range("D32").Value = "'21.52"
$ws.Range("E32").Value = "  +0.82%  "
$ws.Range("D33").Value = "'0.0916"
$ws.Range("E33").Value = "  +0.99%  "
$ws.Range("E34").Value = "  -1.49%  "
$ws.Range("D35").Value = "'4.48"
$ws.Range("E35").Value = "  +14.85%  "
$ws.Range("E36").Value = "  +2.73%  "
$ws.Range("E37").Value = "  +7.32%  "
$ws.Range("D38").Value = "'4.70"
$ws.Range("E38").Value = "  +1.32%  "
$ws.Range("E39").Value = "  +4.66%  "
$ws.Range("E40").Value = "  -4.17%  "
$ws.Range("D41").Value = "'13.74"
$ws.Range("E41").Value = "  -0.92%  "
$ws.Range("E42").Value = "  +1.95%  "
$ws.Range("D43").Value = "'72.29"
$ws.Range("E43").Value = "  -5.94%  "
$ws.Range("E44").Value = "  -0.06%  "
$ws.Range("E45").Value = "  -0.78%  "
$ws.Range("D46").Value = "'5.71"
$ws.Range("E46").Value = "  -8.27%  "
$ws.Range("D47").Value = "'76.86"
$ws.Range("E47").Value = "  +38.84%  "
$ws.Range("D48").Value = "'0.670"
$ws.Range("E48").Value = "  +19.32%  "
$ws.Range("E49").Value = "  +1.67%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.101"
$ws.Range("E50").Value = "  +1.59%  "
$ws.Range("B51").Value = "FraxShare"
$ws.Range("C51").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D51").Value = "'8.57"
$ws.Range("E51").Value = "  -1.22%  "
